$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 29: Mockups til UC03, UC05, UC06, UC10 ---
$ws.Range("A29").Value2 = "Mockups til UC03, UC05, UC06, UC10"
$ws.Range("C29").Value2 = 43969
$ws.Range("D29").Value2 = 0.35416666666666669
$ws.Range("E29").Value2 = 0.5
$ws.Range("F29").Value2 = 0.083333333333333329
$ws.Range("F29").NumberFormat = "h:mm"

# --- Row 30: Samarbejde om DCD og SD ---
$ws.Range("A30").Value2 = "Samarbejde om DCD og SD"
$ws.Range("C30").Value2 = 43969
$ws.Range("D30").Value2 = 0.54166666666666663
$ws.Range("E30").Value2 = 0.58333333333333337
$ws.Range("F30").Value2 = 0.041666666666666664
$ws.Range("F30").NumberFormat = "h:mm"

# --- Row 31: Kundemøde ---
$ws.Range("A31").Value2 = "Kundemøde"
$ws.Range("C31").Value2 = 43969
$ws.Range("D31").Value2 = 0.58333333333333337
$ws.Range("E31").Value2 = 0.625
$ws.Range("F31").Value2 = 0.041666666666666664
$ws.Range("F31").NumberFormat = "h:mm"

# --- Row 32: Samarbejde om DCD og SD ---
$ws.Range("A32").Value2 = "Samarbejde om DCD og SD"
$ws.Range("C32").Value2 = 43969
$ws.Range("D32").Value2 = 0.625
$ws.Range("E32").Value2 = 0.65277777777777779
$ws.Range("F32").Value2 = 0.041666666666666664
$ws.Range("F32").NumberFormat = "h:mm"

# --- Update the view to reflect where the user ended up working ---
[void]$ws.Activate()
[void]$excel.Goto($ws.Range("A10"), $true)
[void]$ws.Range("B33").Select()
